$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '29.485.82'
Set-TextValue $ws.Range('E2') '  -0.93%  '
Set-TextValue $ws.Range('D3') '1.849.10'
Set-TextValue $ws.Range('E3') '  -0.54%  '
Set-TextValue $ws.Range('D4') '0.9990'
Set-TextValue $ws.Range('E4') '  -0.09%  '
Set-TextValue $ws.Range('D5') '241.76'
Set-TextValue $ws.Range('E5') '  -1.02%  '
Set-TextValue $ws.Range('D6') '0.6283'
Set-TextValue $ws.Range('E6') '  -2.15%  '
Set-TextValue $ws.Range('D7') '0.9999'
Set-TextValue $ws.Range('E7') '  -0.06%  '
Set-TextValue $ws.Range('D8') '48.08'
Set-TextValue $ws.Range('E8') '  +0.14%  '
Set-TextValue $ws.Range('D9') '0.07532'
Set-TextValue $ws.Range('D10') '0.2977'
Set-TextValue $ws.Range('E10') '  -0.23%  '
Set-TextValue $ws.Range('D11') '24.36'
Set-TextValue $ws.Range('E11') '  -0.88%  '
Set-TextValue $ws.Range('B12') 'WrappedEther'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D12') '1.957.25'
Set-TextValue $ws.Range('E12') '  +4.97%  '
Set-TextValue $ws.Range('B13') 'TRON'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D13') '0.07718'
Set-TextValue $ws.Range('E13') '  +0.36%  '
Set-TextValue $ws.Range('D14') '5.011'
Set-TextValue $ws.Range('E14') '  -0.78%  '
Set-TextValue $ws.Range('D15') '0.6893'
Set-TextValue $ws.Range('E15') '  -0.32%  '
Set-TextValue $ws.Range('D16') '83.68'
Set-TextValue $ws.Range('E16') '  -0.34%  '
Set-TextValue $ws.Range('D17') '0.000009799'
Set-TextValue $ws.Range('E17') '  -0.63%  '
Set-TextValue $ws.Range('D18') '2.160.94'
Set-TextValue $ws.Range('E18') '  +2.32%  '
Set-TextValue $ws.Range('D19') '6.241'
Set-TextValue $ws.Range('E19') '  +1.90%  '
Set-TextValue $ws.Range('D20') '29.573.36'
Set-TextValue $ws.Range('E20') '  -0.66%  '
Set-TextValue $ws.Range('D21') '233.75'
Set-TextValue $ws.Range('E21') '  -0.95%  '
Set-TextValue $ws.Range('E22') '  -1.24%  '
Set-TextValue $ws.Range('D23') '1.000'
Set-TextValue $ws.Range('E23') '  +0.03%  '
Set-TextValue $ws.Range('D24') '7.620'
Set-TextValue $ws.Range('E24') '  +0.72%  '
Set-TextValue $ws.Range('E25') '  -0.08%  '
Set-TextValue $ws.Range('D26') '154.66'
Set-TextValue $ws.Range('E26') '  -2.47%  '
Set-TextValue $ws.Range('D27') '0.1391'
Set-TextValue $ws.Range('E27') '  -2.21%  '
Set-TextValue $ws.Range('E28') '  -1.27%  '
Set-TextValue $ws.Range('D29') '17.71'
Set-TextValue $ws.Range('E29') '  -1.24%  '
Set-TextValue $ws.Range('D30') '1.479'
Set-TextValue $ws.Range('E30') '  -1.18%  '
Set-TextValue $ws.Range('D31') '0.05847'
Set-TextValue $ws.Range('E31') '  -5.90%  '
Set-TextValue $ws.Range('D32') '1.254'
Set-TextValue $ws.Range('E32') '  -2.80%  '
Set-TextValue $ws.Range('D33') '4.105'
Set-TextValue $ws.Range('E33') '  -1.24%  '
Set-TextValue $ws.Range('D34') '4.036'
Set-TextValue $ws.Range('E34') '  -1.40%  '
Set-TextValue $ws.Range('E35') '  -0.82%  '
Set-TextValue $ws.Range('E36') '  -0.33%  '
Set-TextValue $ws.Range('D37') '0.7207'
Set-TextValue $ws.Range('E37') '  -1.46%  '
Set-TextValue $ws.Range('E38') '  -0.90%  '
Set-TextValue $ws.Range('D39') '1.244.33'
Set-TextValue $ws.Range('E39') '  +2.00%  '
Set-TextValue $ws.Range('D40') '2.798'
Set-TextValue $ws.Range('E40') '  -0.99%  '
Set-TextValue $ws.Range('E41') '  -0.34%  '
Set-TextValue $ws.Range('D42') '0.9071'
Set-TextValue $ws.Range('E42') '  -1.36%  '
Set-TextValue $ws.Range('D43') '6.170'
Set-TextValue $ws.Range('E43') '  -2.25%  '
Set-TextValue $ws.Range('D44') '2.080.33'
Set-TextValue $ws.Range('E44') '  +2.85%  '
Set-TextValue $ws.Range('D45') '0.9996'
Set-TextValue $ws.Range('E45') '  -0.06%  '
Set-TextValue $ws.Range('D46') '102.05'
Set-TextValue $ws.Range('E46') '  +0.09%  '
Set-TextValue $ws.Range('D47') '67.35'
Set-TextValue $ws.Range('E47') '  +0.50%  '
Set-TextValue $ws.Range('D48') '7.327'
Set-TextValue $ws.Range('E48') '  +8.80%  '
Set-TextValue $ws.Range('B49') 'TheSandbox'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D49') '0.4044'
Set-TextValue $ws.Range('E49') '  -0.56%  '
Set-TextValue $ws.Range('B50') 'EnergySwap'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D50') '9.161'
Set-TextValue $ws.Range('E50') '  -0.28%  '
Set-TextValue $ws.Range('B51') 'RenderToken'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D51') '1.707'
Set-TextValue $ws.Range('E51') '  +2.07%  '
